# The exam was re-graded after being trimmed from 3 answer-groups (A/B, D/E, G/H - 56
# questions total) down to 2 groups (A/B always, D/E only for rows 16-18 - 28 questions
# total), and the student's answers were (re-)marked. This script reproduces that
# resulting state on top of the blank "before" marksheet template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the third answer-group (columns G:H) completely - Clear() (not just
#    ClearContents) removes the cells outright so the sheet's used range
#    shrinks from A5:H40 down to A5:E40, same as the "Not Attempt"/"Max"
#    bookkeeping columns D:E beyond the first three questions.
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# 2) Score summary block (rows 10-12). The previous run was a blank template
#    (0 right / 0 wrong / all not-attempted); this is the actual marked
#    result: 9 right, 1 wrong, 18 not attempted out of 28 questions, marked
#    +4 / -1 / 0, i.e. a final score of 35 (9*4 - 1*1) out of 112 (28*4).
# ---------------------------------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 18
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "35/112"

# ---------------------------------------------------------------------------
# 3) Fill in the student's answers in column A for every question that was
#    actually attempted. Style marks right answers green (correctStyle) and
#    wrong answers red (incorrectStyle); attempted cells get both the style
#    and the chosen option text, matching column B (correct answer) when
#    right.
# ---------------------------------------------------------------------------
function Set-Answer($row, $style, $text) {
    $cell = $ws.Range("A$row")
    $cell.Style = $style
    $cell.Value = $text
}

Set-Answer 16 "correctStyle"   "Option A"
Set-Answer 22 "correctStyle"   "Option D"
Set-Answer 23 "correctStyle"   "Option D"
Set-Answer 24 "incorrectStyle" "Option B"
Set-Answer 27 "correctStyle"   "Option A"
Set-Answer 28 "correctStyle"   "Option D"
Set-Answer 30 "correctStyle"   "Option B"
Set-Answer 32 "correctStyle"   "Option C"
Set-Answer 33 "correctStyle"   "Option D"
Set-Answer 38 "correctStyle"   "Option A"
